# Enhance attendance marking and reporting:
# - Split the last date column (T1, "2025-11-22") into a paired check-in/out
#   column (like the existing "_x"/"_y" pair in D1/E1), re-using the existing
#   "Present"/"Total"/"Attendance %" summary columns in between, and append a
#   fresh plain-date column at the very end carrying the original T1 value.
# - Recompute "Total" (expected classes) for every student based on the
#   updated weekly schedule (17 -> 18) and refresh "Attendance %".
# - Automatically mark attendance for recognized students on the new class
#   day: a "pending" mark in the new check-in column, and present/absent in
#   the new check-out column depending on whether the student was recognized.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header -------------------------------------------------------
# Shift T1's current value ("2025-11-22") out to the new last column (Y1),
# copying T1's format (bold / bordered / centered header style) along with
# it, then write the new paired headers into T1 and X1 using the same
# copied format. U1/V1/W1 keep their existing values & format untouched.

$ws.Range("T1").Copy($ws.Range("Y1"))
$ws.Range("T1").Copy($ws.Range("X1"))
$ws.Range("X1").Value = "2025-11-22_y.1"
$ws.Range("T1").Value = "2025-11-22_x.1"

# --- Per-student data rows -----------------------------------------------
# New expected-classes total for every student.
$newTotal = 18

# Roll numbers (by row) that were automatically recognized/marked present
# for the newly added class day.
$recognizedToday = @{
    2  = $false
    3  = $false
    4  = $false
    5  = $false
    6  = $false
    7  = $false
    8  = $false
    9  = $true
    10 = $false
}

foreach ($row in 2..10) {
    $presentCell = $ws.Cells.Item($row, 21)   # column U - Present
    $totalCell   = $ws.Cells.Item($row, 22)   # column V - Total
    $pctCell     = $ws.Cells.Item($row, 23)   # column W - Attendance %
    $checkInCell  = $ws.Cells.Item($row, 24)  # column X - new check-in mark
    $checkOutCell = $ws.Cells.Item($row, 25)  # column Y - new check-out mark

    $present = $presentCell.Value2
    if ($recognizedToday[$row]) {
        $present = $present + 1
    }

    $presentCell.Value = $present
    $totalCell.Value = $newTotal
    $pctCell.Value = [math]::Round(($present / $newTotal) * 100, 1)

    $checkInCell.Value = "⏸️"
    if ($recognizedToday[$row]) {
        $checkOutCell.Value = "✅"
    } else {
        $checkOutCell.Value = "❌"
    }
}
